$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes (C, D, F, H) ---
# The ColumnWidth setter applies Excel's standard character-padding offset
# (~0.8333), so the input is pre-compensated to land on the exact target width
# once saved to OOXML <col width="...">.
$ws.Range("C1").ColumnWidth = 65.16666666666667   # -> stored width 66
$ws.Range("D1").ColumnWidth = 37.166666666666664  # -> stored width 38
$ws.Range("F1").ColumnWidth = 15.166666666666666  # -> stored width 16
$ws.Range("H1").ColumnWidth = 59.166666666666664  # -> stored width 60

# --- Replace data rows 2-10 and append new rows 11-15 ---
# Columns: Id, Title, Country, Premium, Applicants, Duration, Organization
$data = @(
    @('1326953', 'Commodities Support Specialist', 'Panamá, Provincia de Panamá, Panamá', 'No', '2 applicants', '6 - 18 Months', 'NESTRADE S.A, PANAMA BRANCH'),
    @('1326789', 'Content Creation', '4750 Barcelos, Portugal', 'No', '4 applicants', '6 - 18 Months', 'Design Studio'),
    @('1326778', 'Talent Rewards Intern', 'Panamá, Provincia de Panamá, Panamá', 'No', '37 applicants', '6 - 18 Months', 'Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)'),
    @('1326701', 'Software Developer', 'Berlin, Germany', 'No', '18 applicants', '6 - 18 Months', 'code4business Shareholder GmbH'),
    @('1326448', 'TIM Operations Assistant', 'Panamá, Provincia de Panamá, Panamá', 'No', '41 applicants', '6 - 18 Months', 'Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)'),
    @('1326447', 'Controlling CO Intern', 'Panamá, Provincia de Panamá, Panamá', 'No', '26 applicants', '6 - 18 Months', 'Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)'),
    @('1326337', 'Digital Marketing', 'Lisboa, Portugal', 'No', '47 applicants', '6 - 18 Months', 'BGI S.A'),
    @('1326174', 'Accelerate Romania - Digital Marketing & Social Media Assistant', 'Cluj-Napoca, Romania', 'No', '44 applicants', '9 - 12 Weeks', 'Ave Visto'),
    @('1325986', '[Partly Remote] Internal Control Analyst', 'Mexico City, CDMX, Mexico', 'No', '29 applicants', 'Partly Remote', 'Sodexo Mexico'),
    @('1325972', 'Software Engineer', 'Dilovası, Kocaeli, Türkiye', 'No', '30 applicants', '6 - 18 Months', 'AHA TEKNOLOJİ'),
    @('1324560', '[Impact Belo Horiozonte] - Inside Sales Intern', 'Juiz de Fora, MG, Brasil', 'No', '20 applicants', '6 - 18 Months', 'AUE - Juiz de Fora'),
    @('1322487', 'Accelerate Romania - Social Media Manager', 'Cluj-Napoca, Romania', 'No', '54 applicants', '9 - 12 Weeks', 'nclav'),
    @('1317231', 'Guest Relations Executive', 'Katunayake, Sri Lanka', 'No', '13 applicants', '6 - 18 Months', 'Tamarind Tree Garden Resort'),
    @('1314638', 'Accelerate Romania - Business development & Project management', 'Cluj-Napoca, Romania', 'No', '52 applicants', '9 - 12 Weeks', 'ROTSA')
)

$r = 2
foreach ($item in $data) {
    $id = $item[0]
    $ws.Cells.Item($r, 1).Value = "'" + $id
    $ws.Cells.Item($r, 2).Value = "https://aiesec.org/opportunity/global-talent/" + $id
    $ws.Cells.Item($r, 3).Value = $item[1]
    $ws.Cells.Item($r, 4).Value = $item[2]
    $ws.Cells.Item($r, 5).Value = $item[3]
    $ws.Cells.Item($r, 6).Value = $item[4]
    $ws.Cells.Item($r, 7).Value = $item[5]
    $ws.Cells.Item($r, 8).Value = $item[6]
    $r = $r + 1
}

